$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write the date-like label as text (not an auto-converted date serial).
# A plain string assignment of "01-08-2021" gets parsed by Excel as a date,
# so build it as a formula result and paste-special the value through,
# which keeps it a shared-string text cell with the default (no) style.
$ws.Range("Z1").Formula = "=""01-08-2021"""
$ws.Range("Z1").Copy()
$ws.Range("A23").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()
$ws.Range("D23").Value = 12
$ws.Range("E23").Value = 3
$ws.Range("H23").Value = 8.699999999999999
$ws.Range("I23").Value = 3.5
